$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Execute" column (G) with header + values for the two existing rows ---
$ws.Cells.Item(1, 7).Value = "Execute"
$ws.Cells.Item(2, 7).Value = "Yes"

# --- Update existing data row (row 2) with new test values ---
$ws.Cells.Item(2, 2).Value = "qwer"
$ws.Cells.Item(2, 3).Value = "erty"
$ws.Cells.Item(2, 4).Value = "u1,U1,d1"
$ws.Cells.Item(2, 5).Value = "qwer"
$ws.Cells.Item(2, 6).Value = "qwerty"

# --- Add a new data row (row 3) with another test case ---
$ws.Cells.Item(3, 1).Value = "Telecom"
$ws.Cells.Item(3, 2).Value = "asd"
$ws.Cells.Item(3, 3).Value = "uedr"
$ws.Cells.Item(3, 4).Value = "gertf,sddf"
$ws.Cells.Item(3, 5).Value = "asd"
$ws.Cells.Item(3, 6).Value = "asdfg"
$ws.Cells.Item(3, 7).Value = "Yes"

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 18.166666666666668

# --- Update the active selection ---
[void]$ws.Range("H10").Select()

# --- Switch the page to portrait orientation ---
$ws.PageSetup.Orientation = 1
